$wb = $excel.ActiveWorkbook

# --- Add new "Total" worksheet after the last existing sheet ("Tous") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Total"

# --- Header row (annee, sexe, nombre_licencies, Fédération) ---
$ws.Range("A1").Value = "annee"
$ws.Range("B1").Value = "sexe"
$ws.Range("C1").Value = "nombre_licencies"
$ws.Range("D1").Value = "Fédération"

# Re-use the bold/bordered header style already present on the other tabs.
$wb.Worksheets.Item("Tous").Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows ---
$ws.Range("B2").Value = "H"
$ws.Range("C2").Value = 2595
$ws.Range("D2").Value = " Toutes"

$ws.Range("B3").Value = "F"
$ws.Range("C3").Value = 1161
$ws.Range("D3").Value = " Toutes"

# "2021" must land as text (matching the other tabs), not get auto-coerced
# into a number - copy it as a value from existing text cells that already
# hold it (two cells at once so both destination rows get filled).
$wb.Worksheets.Item("F").Range("B2:B3").Copy()
$ws.Range("A2:A3").PasteSpecial(-4163)

# --- Column widths, matching the source diff as closely as this host's
# quantised ColumnWidth setter allows ---
$ws.Columns.Item(1).ColumnWidth = 4.833333333333334
$ws.Columns.Item(2).ColumnWidth = 3.833333333333333
$ws.Columns.Item(3).ColumnWidth = 15.833333333333332
$ws.Columns.Item(4).ColumnWidth = 9.833333333333332
